$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row per the scraped data refresh.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.872.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3662"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07144"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9254"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07686"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.870.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.286"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.392"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008627"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.902.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.928"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.007"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.877"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08814"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.209"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.176"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7457"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.780"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.470"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.085"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01938"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05203"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.952"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5195"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.962"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.150"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4694"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.593"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06036"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  +5.17%  "
